$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect, make the edits, then re-protect
# with the same settings it had before.
$ws.Unprotect()

# Update the confidential disclaimer text with the new "as of" date
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-07 for illustrative purposes only and are subject to change."
# The multi-line text can trigger an automatic row-height bump; AutoFit
# restores the row to its natural (non-custom) height so it round-trips
# the same way it started.
$ws.Rows(11).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5303038492824088
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.2686013116899285
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.05024917889897976
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.09414425349821334
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.02712651462745395
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.02957489200301565
$ws.Range("E7").Value = 0

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0

# Restore sheet protection (contents/objects/scenarios locked; row/column
# formatting allowed), matching the protection that was in place before.
$ws.Protect([Type]::Missing, $true, $true, $true, $true, $true, $true)
